$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$settings.Activate()

$settings.Range("C1").Value = "version"
$settings.Range("C2").Value = 4

$settings.Range("C3").Select()
